$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = 'a'
$ws.Range("B1").Value = 'v'
$ws.Range("A3").Value = 'SILVESTRI 1 Verona'
$ws.Range("B3").Value = 'RADU I. 1 Inter'
$ws.Range("A4").Value = 'SKORUPSKI 1 Bologna'
$ws.Range("B4").Value = 'PANDUR 1 Verona'
$ws.Range("A5").Value = 'MONTIPO'' 1 Benevento'
$ws.Range("B5").Value = 'OSPINA 1 Napoli'
$ws.Range("A7").Value = 'BIRASCHI 1 Genoa'
$ws.Range("B7").Value = 'DIMARCO 1 Verona'
$ws.Range("A8").Value = 'FRABOTTA 1 Juventus'
$ws.Range("B8").Value = 'RUGANI 1 Cagliari'
$ws.Range("A9").Value = 'FOULON 1 Benevento'
$ws.Range("B9").Value = 'GHISLANDI 1 Atalanta'
$ws.Range("A10").Value = 'DARMIAN 1 Inter'
$ws.Range("B10").Value = 'SINGO 1 Torino'
$ws.Range("A11").Value = 'DEPAOLI 1 Benevento'
$ws.Range("B11").Value = 'PELUSO 1 Sassuolo'
$ws.Range("A12").Value = 'DIERCKX 1 Parma'
$ws.Range("B12").Value = 'ARMINI 1 Lazio'
$ws.Range("A13").Value = 'GHOULAM 1 Napoli'
$ws.Range("B13").Value = 'CRISCITO 1 Genoa'
$ws.Range("A14").Value = 'RODRIGO BECAO 1 Udinese'
$ws.Range("B14").Value = 'ROMAGNOLI 1 Milan'
$ws.Range("A16").Value = 'LAZOVIC 1 Verona'
$ws.Range("B16").Value = 'PASTORE 1 Roma'
$ws.Range("A17").Value = 'LUCAS LEIVA 1 Lazio'
$ws.Range("B17").Value = 'SAPONARA 1 Spezia'
$ws.Range("A18").Value = 'JANKTO 1 Sampdoria'
$ws.Range("B18").Value = 'SVANBERG 1 Bologna'
$ws.Range("A19").Value = 'BARAK 1 Verona'
$ws.Range("B19").Value = 'LOBOTKA 1 Napoli'
$ws.Range("A20").Value = 'KULUSEVSKI 1 Juventus'
$ws.Range("B20").Value = 'MANDRAGORA 1 Torino'
$ws.Range("A21").Value = 'NANDEZ 1 Cagliari'
$ws.Range("B21").Value = 'KURTIC 1 Parma'
$ws.Range("A22").Value = 'KUCKA 1 Parma'
$ws.Range("B22").Value = 'POLI 1 Bologna'
$ws.Range("A23").Value = 'JAJALO 1 Udinese'
$ws.Range("B23").Value = 'ARSLAN 1 Udinese'
$ws.Range("A25").Value = 'BRAAF 1 Udinese'
$ws.Range("B25").Value = 'CAPUTO 1 Sassuolo'
$ws.Range("A26").Value = 'SANCHEZ 1 Inter'
$ws.Range("B26").Value = 'HAUGE 1 Milan'
$ws.Range("A27").Value = 'MONCINI 1 Benevento'
$ws.Range("B27").Value = 'JUWARA 1 Bologna'
$ws.Range("A28").Value = 'MESSIAS 1 Crotone'
$ws.Range("B28").Value = 'DRAGUS 1 Crotone'
$ws.Range("A29").Value = 'INGLESE 1 Parma'
$ws.Range("B29").Value = 'LA GUMINA 1 Sampdoria'
$ws.Range("A30").Value = 'ILICIC 1 Atalanta'
$ws.Range("B30").Value = 'DEULOFEU 1 Udinese'
